# Update "想去人数" (F column) figures across all sheets to match the
# latest generated output (gh-pages regeneration at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 896
$ws.Cells.Item(3, 6).Value = 1476
$ws.Cells.Item(4, 6).Value = 1142
$ws.Cells.Item(5, 6).Value = 536
$ws.Cells.Item(6, 6).Value = 231
$ws.Cells.Item(9, 6).Value = 277
$ws.Cells.Item(11, 6).Value = 106
$ws.Cells.Item(13, 6).Value = 166
$ws.Cells.Item(14, 6).Value = 3461
$ws.Cells.Item(15, 6).Value = 16
$ws.Cells.Item(19, 6).Value = 515
$ws.Cells.Item(20, 6).Value = 293
$ws.Cells.Item(22, 6).Value = 117
$ws.Cells.Item(25, 6).Value = 64
$ws.Cells.Item(26, 6).Value = 264
$ws.Cells.Item(27, 6).Value = 974
$ws.Cells.Item(29, 6).Value = 1622
$ws.Cells.Item(30, 6).Value = 355

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 123
$ws.Cells.Item(6, 6).Value = 29
$ws.Cells.Item(7, 6).Value = 243
$ws.Cells.Item(9, 6).Value = 69

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 106

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 896
$ws.Cells.Item(4, 6).Value = 1476
$ws.Cells.Item(5, 6).Value = 1142
$ws.Cells.Item(6, 6).Value = 123
$ws.Cells.Item(8, 6).Value = 106
$ws.Cells.Item(9, 6).Value = 536
$ws.Cells.Item(14, 6).Value = 277
$ws.Cells.Item(16, 6).Value = 106
$ws.Cells.Item(18, 6).Value = 166
$ws.Cells.Item(19, 6).Value = 3462
$ws.Cells.Item(20, 6).Value = 16
$ws.Cells.Item(25, 6).Value = 515
$ws.Cells.Item(26, 6).Value = 293
$ws.Cells.Item(28, 6).Value = 29
$ws.Cells.Item(29, 6).Value = 117
$ws.Cells.Item(31, 6).Value = 243
$ws.Cells.Item(33, 6).Value = 69
$ws.Cells.Item(38, 6).Value = 64
$ws.Cells.Item(39, 6).Value = 264
$ws.Cells.Item(40, 6).Value = 974
$ws.Cells.Item(42, 6).Value = 1622
$ws.Cells.Item(43, 6).Value = 355
